$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = 6.4
$ws.Range("A12").Value = -21.448
$ws.Range("C14").Value = -12.697
$ws.Range("C19").Value = -12.177
$ws.Range("B23").Value = 7.513000000000001
$ws.Range("C24").Value = -12.641
$ws.Range("A27").Value = -21.652
$ws.Range("B28").Value = 5.172
$ws.Range("A32").Value = -20.615
$ws.Range("B32").Value = 7.537999999999999
$ws.Range("B34").Value = 6.572
$ws.Range("A36").Value = -20.512
$ws.Range("A38").Value = -20.474
$ws.Range("C38").Value = -11.626
$ws.Range("C41").Value = -11.991
$ws.Range("B42").Value = 7.761
$ws.Range("A46").Value = -21.651
$ws.Range("B49").Value = 6.237
$ws.Range("C52").Value = -11.632
$ws.Range("A54").Value = -20.82
$ws.Range("B54").Value = 6.215999999999999
$ws.Range("A55").Value = -22.184
$ws.Range("A56").Value = -21.461
$ws.Range("A67").Value = -21.418
$ws.Range("A69").Value = -21.323
$ws.Range("A72").Value = -21.194
$ws.Range("C72").Value = -12.648
$ws.Range("B78").Value = 7.329000000000001
$ws.Range("C78").Value = -11.833
$ws.Range("B80").Value = 7.398999999999999
$ws.Range("A83").Value = -21.294
$ws.Range("C83").Value = -13.419
$ws.Range("C85").Value = -12.335
$ws.Range("A86").Value = -21.429
$ws.Range("C86").Value = -13.349
$ws.Range("C90").Value = -10.589
$ws.Range("A91").Value = -20.905
$ws.Range("A93").Value = -21.338
$ws.Range("C96").Value = -10.38
$ws.Range("B97").Value = 5.401999999999999
$ws.Range("A99").Value = -20.861
$ws.Range("B99").Value = 6.287999999999999
$ws.Range("B101").Value = 5.661
$ws.Range("C103").Value = -12.411
$ws.Range("A104").Value = -21.437
